$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" date) for rows 2-8 from 45183 to 45184.
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45184
}
